# Refresh crypto price (D) / volume (E) figures per the scraper run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cell -> new text value, exactly as scraped (preserves the source
# formatting: trailing zeros, dot-grouped thousands, padded "  +x.xx%  ").
$updates = [ordered]@{
    'D2' = '27.831.55'
    'E2' = '  -0.87%  '
    'D3' = '1.902.26'
    'E3' = '  -0.38%  '
    'E4' = '  -0.65%  '
    'D5' = '312.78'
    'E5' = '  -1.23%  '
    'E6' = '  -0.62%  '
    'D7' = '0.5009'
    'D8' = '0.3809'
    'D9' = '0.07291'
    'E9' = '  -0.89%  '
    'D10' = '0.9092'
    'E10' = '  -2.61%  '
    'D11' = '20.87'
    'E11' = '  +0.70%  '
    'D12' = '0.07643'
    'E12' = '  -2.31%  '
    'D13' = '1.922.72'
    'E13' = '  +0.21%  '
    'D14' = '5.482'
    'E14' = '  -0.11%  '
    'D15' = '6.612'
    'E15' = '  -0.07%  '
    'D16' = '91.28'
    'E16' = '  +0.17%  '
    'E17' = '  -0.77%  '
    'D18' = '0.000008715'
    'E18' = '  -1.12%  '
    'D19' = '1.000'
    'E19' = '  -0.66%  '
    'D20' = '27.870.42'
    'E20' = '  -0.86%  '
    'D21' = '14.50'
    'E21' = '  -2.05%  '
    'D22' = '5.148'
    'E22' = '  -0.07%  '
    'D23' = '10.82'
    'E23' = '  -0.39%  '
    'D24' = '154.45'
    'E24' = '  -1.42%  '
    'D25' = '1.861'
    'E25' = '  -2.71%  '
    'D26' = '2.232'
    'E26' = '  +6.90%  '
    'D27' = '18.37'
    'E27' = '  -0.89%  '
    'D28' = '115.24'
    'E28' = '  -0.81%  '
    'D29' = '4.931'
    'E29' = '  -0.31%  '
    'D30' = '0.08973'
    'E30' = '  +0.79%  '
    'D31' = '3.205'
    'E31' = '  -4.87%  '
    'D32' = '1.236'
    'E32' = '  -0.46%  '
    'D33' = '0.7704'
    'E33' = '  +0.57%  '
    'D34' = '4.640'
    'E34' = '  -0.76%  '
    'D35' = '0.02060'
    'E35' = '  +0.88%  '
    'D36' = '2.562'
    'E36' = '  -1.19%  '
    'D37' = '1.100'
    'E37' = '  +0.23%  '
    'D38' = '0.5535'
    'E38' = '  +0.85%  '
    'D39' = '3.013'
    'E39' = '  +0.61%  '
    'D40' = '0.05273'
    'E40' = '  -0.40%  '
    'D41' = '6.984'
    'E41' = '  -0.63%  '
    'D42' = '8.548'
    'E42' = '  +1.57%  '
    'D43' = '0.1523'
    'E43' = '  +0.12%  '
    'D44' = '111.09'
    'E44' = '  +4.02%  '
    'D45' = '10.60'
    'E45' = '  -1.05%  '
    'D46' = '0.4798'
    'E46' = '  -0.73%  '
    'E47' = '  -0.68%  '
    'D48' = '1.639'
    'E48' = '  -0.86%  '
    'D49' = '67.31'
    'E49' = '  -1.58%  '
    'E50' = '  -0.31%  '
    'D51' = '0.9003'
    'E51' = '  -0.32%  '
}

foreach ($cell in $updates.Keys) {
    $value = $updates[$cell]
    $range = $ws.Range($cell)
    if ($cell -match "^D" -and $value -match "^[0-9]+\.[0-9]+$") {
        # Plain single-dot decimals (e.g. "14.50", "0.02060") would otherwise be
        # coerced to a Double by the Value setter and lose exact text formatting,
        # so force text storage; restore the default style afterwards so no new
        # cell style is left applied (matches the unstyled source cells).
        $range.NumberFormat = "@"
        $range.Value = $value
        $range.Style = "Normal"
    } else {
        $range.Value = $value
    }
}
